$d = $word.ActiveDocument

# --- 1) First paragraph: pad trailing text with two spaces, then append a
#        red-colored "(This is a change - Version for branch alternate)"
#        split across three runs, reproducing the authored edit. ---
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$p1 = $d.Paragraphs.Item(1)
$ins = $p1.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)   # collapse to just before the paragraph mark

$ins.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$ins.Font.Color = 192            # RGB(192,0,0) -> wdColor BGR long -> w:val="C00000"
$ins.Collapse(0)                 # wdCollapseEnd

$ins.InsertAfter("rsion for branch alternate")
$ins.Font.Color = 192
$ins.Collapse(0)

$ins.InsertAfter(")")
$ins.Font.Color = 192

# --- 2) The blank paragraph right after "It will be treated..." gets new
#        shading/formatting (Calibri, bold, color 202122, fill F9F9F9) in
#        place of the old Menlo-flavoured blank line. ---
$blank = $d.Paragraphs.Item(3)
$blankXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
              '<w:pPr>' + `
                '<w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/>' + `
                '<w:rPr>' + `
                  '<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>' + `
                  '<w:b/>' + `
                  '<w:bCs/>' + `
                  '<w:color w:val="202122"/>' + `
                '</w:rPr>' + `
              '</w:pPr>' + `
            '</w:p>'
$blank.Range.InsertXML($blankXml) | Out-Null

# --- 3) Drop the trailing "ank God almighty, we are free at last." NormalWeb
#        paragraph, leaving a single empty paragraph before the sectPr. ---
$lastIdx = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($lastIdx)
$last.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>') | Out-Null
